# Bug fix: "header text with link resources" was being omitted because
# internal `_h2o_keep_element` placeholder runs on certain structural
# paragraphs were left as literal text instead of being collapsed down
# to a single blank space once the surrounding content was resolved.
#
# Rule (matches the fix): every "Node End" paragraph collapses to " ".
# A "Head End" paragraph collapses to " " only when it is immediately
# followed by a "Body Text" paragraph (i.e. it terminates the heading of
# a resource that has inline body content). "Head End" paragraphs that
# close out the document's top-level head, or that close a resource
# which is itself a link (no following Body Text), keep their
# `_h2o_keep_element` marker untouched.

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($text -ne "_h2o_keep_element") {
        continue
    }

    $shouldBlank = $false

    if ($styleName -eq "Node End") {
        $shouldBlank = $true
    }
    elseif ($styleName -eq "Head End") {
        if ($i -lt $count) {
            $nextStyle = $d.Paragraphs.Item($i + 1).Style.NameLocal
            if ($nextStyle -eq "Body Text") {
                $shouldBlank = $true
            }
        }
    }

    if ($shouldBlank) {
        $p.Range.Text = " "
    }
}
